# Updates cryptos list data (prices in column D, 1h volume % in column E,
# plus a row-content shift in B/C/D/E for rows 45-51) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.950.35"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "'1.621.01"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'212.57"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'18.41"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D12").Value = "'1.846.22"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "'1.624.66"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "'4.14"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "'25.959.00"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "'61.71"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "'192.04"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'9.53"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "'6.03"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'0.133"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'143.79"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "'1.71"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").Value = "'15.25"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "'0.0479"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "'1.126.07"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'0.844"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'0.515"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "'97.68"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "'1.757.65"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "'0.758"
$ws.Range("E43").Value = "  -3.61%  "
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.52"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'54.12"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0517"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.410"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.48"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0925"
$ws.Range("E51").Value = "  -0.45%  "
